# "added some rf components"
# Rename the sheet (ADC -> rf) and append five new RF-component rows
# that were sourced / typed column-B-first, then A, then C..G (matching
# the order the parts datasheets were transcribed in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "rf"

# Part Number | Size | Part Number Ru | Part Number En | Library Ref | Footprint Ref | Manufacturer
$newRows = @(
    @("ADE-1+",              "CD636",        "ADE-1+ CD636",               "ADE-1+",              "ADE-1+",              "CD636",        "Mimi-Circuits"),
    @("MAAM-011206",         "6-TDFN",       "MAAM-011206 6-TDFN",         "MAAM-011206",         "MAAM-011206",         "6-TDFN",       "MACOM"),
    @("MADL-011023-14150T",  "6-TDFN",       "MADL-011023-14150T 6-TDFN",  "MADL-011023-14150T",  "MADL-011023-14150T",  "6-TDFN",       "MACOM"),
    @("MTX2-73+",            "12-QFN (3Х3)", "MTX2-73+ 12-QFN (3Х3)",      "MTX2-73+",            "MTX2-73+",            "12-QFN (3Х3)", "Mimi-Circuits"),
    @("NCS2-622+",           "PL-264",       "NCS2-622+ PL-264",           "NCS2-622+",           "NCS2-622+",           "PL-264",       "Mimi-Circuits")
)

# Columns get filled B, A, C, D, E, F, G per row (0-based indexes into each row array)
$fillOrder = @(1, 0, 2, 3, 4, 5, 6)

$rowIndex = 12
foreach ($row in $newRows) {
    foreach ($colOffset in $fillOrder) {
        $ws.Cells.Item($rowIndex, $colOffset + 1).Value = $row[$colOffset]
    }
    $rowIndex++
}

# Column widths widened to fit the new, longer part numbers
$ws.Columns.Item(1).ColumnWidth = 21
$ws.Columns.Item(4).ColumnWidth = 21.6
$ws.Columns.Item(5).ColumnWidth = 20.15

# Leave the selection where the author last clicked
[void]$ws.Range("D9").Select()
